# Inserts a new weekly record as row 3 (shifting all existing data rows
# 3-48 down to rows 4-49) in the Jengibre / Terminal La Palmera de La
# Serena price sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing rows 3:48 down to 4:49, duplicating row 3's formatting
# (matches the style shift seen for column D in the target workbook).
$ws.Rows("3:3").Insert()

# Populate the newly inserted row 3 with the latest week's observation.
$ws.Cells.Item(3, 1).Value2  = 8
$ws.Cells.Item(3, 2).Value2  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(3, 3).Value2  = "Coquimbo"
$ws.Cells.Item(3, 4).Value2  = 44750
$ws.Cells.Item(3, 5).Value2  = 4
$ws.Cells.Item(3, 6).Value2  = 100114007
$ws.Cells.Item(3, 7).Value2  = "Jengibre"
$ws.Cells.Item(3, 8).Value2  = "Sin especificar"
$ws.Cells.Item(3, 9).Value2  = "Primera"
$ws.Cells.Item(3, 10).Value2 = 480
$ws.Cells.Item(3, 11).Value2 = 15000
$ws.Cells.Item(3, 12).Value2 = 16000
$ws.Cells.Item(3, 13).Value2 = 15500
$ws.Cells.Item(3, 14).Value2 = "`$/caja 13 kilos"
$ws.Cells.Item(3, 15).Value2 = "Perú"
$ws.Cells.Item(3, 16).Value2 = 1192
$ws.Cells.Item(3, 17).Value2 = 13
$ws.Cells.Item(3, 18).Value2 = "Hortaliza"
